$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column letter, new value
$changes = @(
    @{Row=2; Col='D'; Value='66.811.86'},
    @{Row=2; Col='E'; Value='  -2.58%  '},
    @{Row=3; Col='D'; Value='2.455.58'},
    @{Row=3; Col='E'; Value='  -3.87%  '},
    @{Row=4; Col='E'; Value='  -0.06%  '},
    @{Row=5; Col='D'; Value='580.53'},
    @{Row=5; Col='E'; Value='  -2.26%  '},
    @{Row=6; Col='D'; Value='167.77'},
    @{Row=6; Col='E'; Value='  -4.86%  '},
    @{Row=8; Col='E'; Value='  -3.20%  '},
    @{Row=9; Col='D'; Value='2.455.69'},
    @{Row=9; Col='E'; Value='  -3.96%  '},
    @{Row=10; Col='E'; Value='  -4.09%  '},
    @{Row=11; Col='E'; Value='  -0.94%  '},
    @{Row=12; Col='D'; Value='4.87'},
    @{Row=12; Col='E'; Value='  -3.24%  '},
    @{Row=13; Col='E'; Value='  -5.69%  '},
    @{Row=14; Col='D'; Value='2.905.23'},
    @{Row=14; Col='E'; Value='  -1.75%  '},
    @{Row=15; Col='D'; Value='25.21'},
    @{Row=15; Col='E'; Value='  -5.40%  '},
    @{Row=16; Col='D'; Value='66.655.36'},
    @{Row=16; Col='E'; Value='  -2.51%  '},
    @{Row=17; Col='E'; Value='  -5.91%  '},
    @{Row=18; Col='D'; Value='2.469.30'},
    @{Row=18; Col='E'; Value='  -8.10%  '},
    @{Row=19; Col='E'; Value='  -8.48%  '},
    @{Row=20; Col='D'; Value='7.39'},
    @{Row=20; Col='E'; Value='  -8.20%  '},
    @{Row=21; Col='D'; Value='349.28'},
    @{Row=21; Col='E'; Value='  -6.11%  '},
    @{Row=22; Col='D'; Value='4.01'},
    @{Row=22; Col='E'; Value='  -4.12%  '},
    @{Row=23; Col='D'; Value='1.00'},
    @{Row=23; Col='E'; Value='  +0.04%  '},
    @{Row=24; Col='D'; Value='68.57'},
    @{Row=24; Col='E'; Value='  -4.82%  '},
    @{Row=25; Col='E'; Value='  -8.90%  '},
    @{Row=26; Col='D'; Value='1.80'},
    @{Row=26; Col='E'; Value='  -6.22%  '},
    @{Row=27; Col='D'; Value='9.05'},
    @{Row=27; Col='E'; Value='  -9.37%  '},
    @{Row=28; Col='D'; Value='0.998'},
    @{Row=28; Col='E'; Value='  -43.92%  '},
    @{Row=29; Col='D'; Value='2.583.88'},
    @{Row=29; Col='E'; Value='  -3.26%  '},
    @{Row=30; Col='D'; Value='0.0₃0889'},
    @{Row=30; Col='E'; Value='  -8.81%  '},
    @{Row=31; Col='D'; Value='505.56'},
    @{Row=31; Col='E'; Value='  -6.20%  '},
    @{Row=32; Col='E'; Value='  -8.79%  '},
    @{Row=33; Col='E'; Value='  -6.62%  '},
    @{Row=34; Col='E'; Value='  -8.65%  '},
    @{Row=35; Col='E'; Value='  -0.11%  '},
    @{Row=36; Col='D'; Value='158.11'},
    @{Row=36; Col='E'; Value='  -1.41%  '},
    @{Row=37; Col='E'; Value='  -12.76%  '},
    @{Row=38; Col='D'; Value='18.63'},
    @{Row=38; Col='E'; Value='  -0.12%  '},
    @{Row=39; Col='E'; Value='  -5.98%  '},
    @{Row=40; Col='E'; Value='  -9.53%  '},
    @{Row=41; Col='E'; Value='  +0.28%  '},
    @{Row=42; Col='D'; Value='1.66'},
    @{Row=42; Col='E'; Value='  -7.34%  '},
    @{Row=43; Col='D'; Value='4.75'},
    @{Row=43; Col='E'; Value='  -8.38%  '},
    @{Row=44; Col='D'; Value='0.322'},
    @{Row=44; Col='E'; Value='  -8.28%  '},
    @{Row=45; Col='E'; Value='  -8.68%  '},
    @{Row=46; Col='D'; Value='38.52'},
    @{Row=46; Col='E'; Value='  -2.46%  '},
    @{Row=47; Col='D'; Value='140.07'},
    @{Row=47; Col='E'; Value='  -6.13%  '},
    @{Row=48; Col='E'; Value='  -8.75%  '},
    @{Row=49; Col='D'; Value='0.505'},
    @{Row=49; Col='E'; Value='  -9.16%  '},
    @{Row=50; Col='E'; Value='  -10.64%  '},
    @{Row=51; Col='D'; Value='0.0726'},
    @{Row=51; Col='E'; Value='  -2.96%  '}
)

foreach ($chg in $changes) {
    $cellRef = "$($chg.Col)$($chg.Row)"
    $cell = $ws.Range($cellRef)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $chg.Value
    $cell.Style = $origStyle
}

